$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level metadata (file path / revision / window position) ---

# --- Cell value updates (Experiment grids) ---
$ws.Range("C6").Value = 3.26
$ws.Range("D6").Value = 6.52
$ws.Range("E6").Value = 9.7799999999999994
$ws.Range("F6").Value = 13.04
$ws.Range("G6").Value = 16.3
$ws.Range("H6").Value = 19.559999999999999
$ws.Range("I6").Value = 22.82
$ws.Range("J6").Value = 26.08
$ws.Range("K6").Value = 29.34
$ws.Range("L6").Value = 32.6
$ws.Range("C8").Value = 4.5599999999999996
$ws.Range("D8").Value = 9.1199999999999992
$ws.Range("E8").Value = 13.68
$ws.Range("F8").Value = 18.239999999999998
$ws.Range("G8").Value = 22.8
$ws.Range("H8").Value = 27.36
$ws.Range("I8").Value = 31.92
$ws.Range("J8").Value = 36.479999999999997
$ws.Range("K8").Value = 41.04
$ws.Range("L8").Value = 45.6
$ws.Range("C15").Value = 16.3
$ws.Range("D15").Value = 32.6
$ws.Range("E15").Value = 48.9
$ws.Range("F15").Value = 65.2
$ws.Range("G15").Value = 81.5
$ws.Range("H15").Value = 97.8
$ws.Range("I15").Value = 114.1
$ws.Range("J15").Value = 130.4
$ws.Range("K15").Value = 146.69999999999999
$ws.Range("L15").Value = 163
$ws.Range("C16").Value = 22.8
$ws.Range("D16").Value = 45.6
$ws.Range("E16").Value = 68.400000000000006
$ws.Range("F16").Value = 91.2
$ws.Range("G16").Value = 114
$ws.Range("H16").Value = 136.80000000000001
$ws.Range("I16").Value = 159.6
$ws.Range("J16").Value = 182.4
$ws.Range("K16").Value = 205.2
$ws.Range("L16").Value = 228
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 20
$ws.Range("E20").Value = 30
$ws.Range("F20").Value = 40
$ws.Range("G20").Value = 50
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 70
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = 90
$ws.Range("L20").Value = 100
$ws.Range("C21").Value = 4.827
$ws.Range("D21").Value = 7.4379999999999997
$ws.Range("E21").Value = 10.041
$ws.Range("F21").Value = 12.638999999999999
$ws.Range("G21").Value = 15.250999999999999
$ws.Range("H21").Value = 17.841999999999999
$ws.Range("I21").Value = 20.440999999999999
$ws.Range("J21").Value = 23.062000000000001
$ws.Range("K21").Value = 25.645
$ws.Range("L21").Value = 28.239000000000001
$ws.Range("C22").Value = 5.8579999999999997
$ws.Range("D22").Value = 9.4320000000000004
$ws.Range("E22").Value = 13.016
$ws.Range("F22").Value = 16.623999999999999
$ws.Range("G22").Value = 20.286999999999999
$ws.Range("H22").Value = 23.876999999999999
$ws.Range("I22").Value = 27.398
$ws.Range("J22").Value = 31.212
$ws.Range("K22").Value = 34.460999999999999
$ws.Range("L22").Value = 38.731999999999999
$ws.Range("C23").Value = 4.093
$ws.Range("D23").Value = 4.9560000000000004
$ws.Range("E23").Value = 5.6349999999999998
$ws.Range("F23").Value = 6.4729999999999999
$ws.Range("G23").Value = 7.3049999999999997
$ws.Range("H23").Value = 8.1649999999999991
$ws.Range("I23").Value = 9.0950000000000006
$ws.Range("J23").Value = 9.8420000000000005
$ws.Range("K23").Value = 10.675000000000001
$ws.Range("L23").Value = 11.507999999999999
$ws.Range("C24").Value = 16.3
$ws.Range("D24").Value = 32.6
$ws.Range("E24").Value = 48.9
$ws.Range("F24").Value = 65.2
$ws.Range("G24").Value = 81.5
$ws.Range("H24").Value = 97.8
$ws.Range("I24").Value = 114.1
$ws.Range("J24").Value = 130.4
$ws.Range("K24").Value = 146.69999999999999
$ws.Range("L24").Value = 163
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 68
$ws.Range("E25").Value = 102
$ws.Range("F25").Value = 136
$ws.Range("G25").Value = 170
$ws.Range("H25").Value = 204
$ws.Range("I25").Value = 238
$ws.Range("J25").Value = 272
$ws.Range("K25").Value = 306
$ws.Range("L25").Value = 340
$ws.Range("C26").Value = 22.8
$ws.Range("D26").Value = 45.6
$ws.Range("E26").Value = 68.400000000000006
$ws.Range("F26").Value = 91.2
$ws.Range("G26").Value = 114
$ws.Range("H26").Value = 136.80000000000001
$ws.Range("I26").Value = 159.6
$ws.Range("J26").Value = 182.4
$ws.Range("K26").Value = 205.2
$ws.Range("C38").Value = 6
$ws.Range("D38").Value = 12
$ws.Range("E38").Value = 18
$ws.Range("F38").Value = 24
$ws.Range("G38").Value = 30
$ws.Range("H38").Value = 36
$ws.Range("I38").Value = 42
$ws.Range("J38").Value = 48
$ws.Range("K38").Value = 54
$ws.Range("L38").Value = 60
$ws.Range("C39").Value = 0.83099999999999996
$ws.Range("D39").Value = 1.6890000000000001
$ws.Range("E39").Value = 2.5019999999999998
$ws.Range("F39").Value = 3.3239999999999998
$ws.Range("G39").Value = 4.0679999999999996
$ws.Range("H39").Value = 4.8849999999999998
$ws.Range("I39").Value = 5.8140000000000001
$ws.Range("J39").Value = 6.6459999999999999
$ws.Range("K39").Value = 7.4359999999999999
$ws.Range("L39").Value = 8.4090000000000007
$ws.Range("C40").Value = 1.419
$ws.Range("D40").Value = 2.859
$ws.Range("E40").Value = 4.1280000000000001
$ws.Range("F40").Value = 6.0179999999999998
$ws.Range("G40").Value = 7.1550000000000002
$ws.Range("H40").Value = 8.5920000000000005
$ws.Range("I40").Value = 9.9510000000000005
$ws.Range("J40").Value = 11.778
$ws.Range("K40").Value = 12.946999999999999
$ws.Range("L40").Value = 14.045999999999999
$ws.Range("C41").Value = 0.64500000000000002
$ws.Range("D41").Value = 1.276
$ws.Range("E41").Value = 1.913
$ws.Range("F41").Value = 2.504
$ws.Range("G41").Value = 3.2530000000000001
$ws.Range("H41").Value = 3.8370000000000002
$ws.Range("I41").Value = 4.5419999999999998
$ws.Range("J41").Value = 5.319
$ws.Range("K41").Value = 5.8289999999999997
$ws.Range("L41").Value = 6.282
$ws.Range("C42").Value = 7.92
$ws.Range("D42").Value = 15.84
$ws.Range("E42").Value = 23.76
$ws.Range("F42").Value = 31.68
$ws.Range("G42").Value = 39.6
$ws.Range("H42").Value = 47.52
$ws.Range("I42").Value = 55.44
$ws.Range("J42").Value = 63.36
$ws.Range("K42").Value = 71.28
$ws.Range("L42").Value = 79.2
$ws.Range("C43").Value = 15.6
$ws.Range("D43").Value = 31.2
$ws.Range("E43").Value = 46.8
$ws.Range("F43").Value = 62.4
$ws.Range("G43").Value = 78
$ws.Range("H43").Value = 93.6
$ws.Range("I43").Value = 109.2
$ws.Range("J43").Value = 124.8
$ws.Range("K43").Value = 140.4
$ws.Range("L43").Value = 156
$ws.Range("C44").Value = 6
$ws.Range("D44").Value = 12
$ws.Range("E44").Value = 18
$ws.Range("F44").Value = 24
$ws.Range("G44").Value = 30
$ws.Range("H44").Value = 36
$ws.Range("I44").Value = 42
$ws.Range("J44").Value = 48
$ws.Range("K44").Value = 54
$ws.Range("L44").Value = 60
$ws.Range("E50").Value = 1
$ws.Range("F50").Value = 1

# --- Clear now-empty cells (row 48-50 trimmed ranges) ---
$ws.Range("G48,H48,I48,J48,K48,L48,M48,N48,G49,H49,I49,J49,K49,L49,M49,N49,G50,H50,I50,J50,K50,L50,M50,N50").ClearContents()

# --- Sheet view (scroll position & selection) ---
$ws.Range("F50").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
